$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (no numeric auto-conversion / no trailing-zero loss)
# for the Price (D) cells whose new values look numeric, matching the
# original inline-string cell type.
$ws.Range("D4:D8").NumberFormat = "@"
$ws.Range("D10:D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19:D29").NumberFormat = "@"
$ws.Range("D31:D34").NumberFormat = "@"
$ws.Range("D36:D44").NumberFormat = "@"
$ws.Range("D46:D51").NumberFormat = "@"

# Apply the updated cell values row by row, in sheet order.

# Row 2
$ws.Range("D2").Value = "66.905.69"
$ws.Range("E2").Value = "  -0.53%  "

# Row 3
$ws.Range("D3").Value = "3.511.11"
$ws.Range("E3").Value = "  +0.79%  "

# Row 4
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "584.14"
$ws.Range("E5").Value = "  -1.60%  "

# Row 6
$ws.Range("D6").Value = "177.72"
$ws.Range("E6").Value = "  -0.44%  "

# Row 7
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.08%  "

# Row 8
$ws.Range("D8").Value = "0.600"
$ws.Range("E8").Value = "  +0.92%  "

# Row 9
$ws.Range("D9").Value = "3.508.86"
$ws.Range("E9").Value = "  +0.63%  "

# Row 10
$ws.Range("D10").Value = "0.136"
$ws.Range("E10").Value = "  -0.88%  "

# Row 11
$ws.Range("D11").Value = "6.92"
$ws.Range("E11").Value = "  -2.22%  "

# Row 12
$ws.Range("D12").Value = "0.421"
$ws.Range("E12").Value = "  -3.20%  "

# Row 13
$ws.Range("D13").Value = "4.104.47"
$ws.Range("E13").Value = "  +0.41%  "

# Row 14
$ws.Range("D14").Value = "30.64"
$ws.Range("E14").Value = "  -4.34%  "

# Row 15
$ws.Range("E15").Value = "  -2.72%  "

# Row 16
$ws.Range("D16").Value = "66.860.57"
$ws.Range("E16").Value = "  -0.68%  "

# Row 17
$ws.Range("D17").Value = "0.0000174"
$ws.Range("E17").Value = "  -1.70%  "

# Row 18
$ws.Range("D18").Value = "3.480.93"
$ws.Range("E18").Value = "  -0.03%  "

# Row 19
$ws.Range("D19").Value = "6.07"
$ws.Range("E19").Value = "  -3.08%  "

# Row 20
$ws.Range("D20").Value = "14.02"
$ws.Range("E20").Value = "  -1.99%  "

# Row 21
$ws.Range("D21").Value = "381.15"
$ws.Range("E21").Value = "  -2.00%  "

# Row 22
$ws.Range("D22").Value = "7.85"
$ws.Range("E22").Value = "  -1.89%  "

# Row 23
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.17%  "

# Row 24
$ws.Range("D24").Value = "5.75"
$ws.Range("E24").Value = "  +0.55%  "

# Row 25
$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").Value = "0.533"
$ws.Range("E25").Value = "  -0.64%  "

# Row 26
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").Value = "71.62"
$ws.Range("E26").Value = "  -3.42%  "

# Row 27
$ws.Range("D27").Value = "0.0000121"
$ws.Range("E27").Value = "  -0.18%  "

# Row 28
$ws.Range("D28").Value = "9.93"
$ws.Range("E28").Value = "  -4.29%  "

# Row 29
$ws.Range("D29").Value = "0.174"
$ws.Range("E29").Value = "  -0.33%  "

# Row 30
$ws.Range("E30").Value = "  -0.03%  "

# Row 31
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "5.99"
$ws.Range("E31").Value = "  -2.80%  "

# Row 32
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "24.41"
$ws.Range("E32").Value = "  +3.65%  "

# Row 33
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "1.37"
$ws.Range("E33").Value = "  -3.58%  "

# Row 34
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").Value = "2.01"
$ws.Range("E34").Value = "  -2.53%  "

# Row 35
$ws.Range("E35").Value = "  -0.07%  "

# Row 36
$ws.Range("D36").Value = "7.18"
$ws.Range("E36").Value = "  -2.66%  "

# Row 37
$ws.Range("D37").Value = "1.54"
$ws.Range("E37").Value = "  -3.18%  "

# Row 38
$ws.Range("D38").Value = "158.81"
$ws.Range("E38").Value = "  -3.18%  "

# Row 39
$ws.Range("D39").Value = "0.887"
$ws.Range("E39").Value = "  +1.71%  "

# Row 40
$ws.Range("D40").Value = "28.20"
$ws.Range("E40").Value = "  +7.39%  "

# Row 41
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "1.81"
$ws.Range("E41").Value = "  -3.64%  "

# Row 42
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "2.66"
$ws.Range("E42").Value = "  -3.04%  "

# Row 43
$ws.Range("D43").Value = "6.55"
$ws.Range("E43").Value = "  -3.99%  "

# Row 44
$ws.Range("D44").Value = "4.51"
$ws.Range("E44").Value = "  -3.21%  "

# Row 45
$ws.Range("D45").Value = "2.718.53"
$ws.Range("E45").Value = "  -4.34%  "

# Row 46
$ws.Range("D46").Value = "0.0704"
$ws.Range("E46").Value = "  -2.62%  "

# Row 47
$ws.Range("D47").Value = "25.55"
$ws.Range("E47").Value = "  -6.18%  "

# Row 48
$ws.Range("D48").Value = "40.45"
$ws.Range("E48").Value = "  -2.89%  "

# Row 49
$ws.Range("D49").Value = "0.0298"
$ws.Range("E49").Value = "  -0.80%  "

# Row 50
$ws.Range("D50").Value = "326.30"
$ws.Range("E50").Value = "  -2.80%  "

# Row 51
$ws.Range("D51").Value = "1.03"
$ws.Range("E51").Value = "  -2.24%  "
